# chore(runtime): publish files + archive (2025-11-29 11:05:00)
#
# Adds the 2025-11-28 KHL match day to the workbook:
#   - Matches_SOG : 5 new match rows (330-334)
#   - Shots_HA    : refreshed as_of_utc + home/away shots-on-goal aggregates
#   - Shots_Summary: refreshed as_of_utc + total shots-on-goal aggregates
#   - Meta_ext    : refreshed as_of_utc + build_version

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Matches_SOG - append the 5 games played on 2025-11-28
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @{row=330; uid="897827"; date="2025-11-28T17:00:00"; home="Трактор";     away="ХК Сочи";     sogHome=33; sogAway=32},
    @{row=331; uid="897823"; date="2025-11-28T19:30:00"; home="Ак Барс";     away="СКА";          sogHome=39; sogAway=29},
    @{row=332; uid="897824"; date="2025-11-28T19:30:00"; home="Торпедо";     away="Динамо Мн";    sogHome=37; sogAway=31},
    @{row=333; uid="897825"; date="2025-11-28T19:30:00"; home="Северсталь"; away="Локомотив";    sogHome=18; sogAway=32},
    @{row=334; uid="897826"; date="2025-11-28T19:00:00"; home="Нефтехимик"; away="Драконы";      sogHome=39; sogAway=50}
)

foreach ($m in $newMatches) {
    # uid column must be stored as text, like the rest of column A
    $ws1.Cells.Item($m.row, 1).NumberFormat = "@"
    $ws1.Cells.Item($m.row, 1).Value = $m.uid

    $ws1.Cells.Item($m.row, 2).Value = $m.date
    $ws1.Cells.Item($m.row, 3).Value = $m.home
    $ws1.Cells.Item($m.row, 4).Value = $m.away
    $ws1.Cells.Item($m.row, 5).Value = $m.sogHome
    $ws1.Cells.Item($m.row, 6).Value = $m.sogAway
    $ws1.Cells.Item($m.row, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------
# 2. Shots_HA - per-team home/away shots-on-goal, as of 2025-11-28
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Shots_HA")
$asOfUtc = "2025-11-28T19:30:00Z"

$shotsHa = @(
    @{row=2; E=14; F=15; G=449; H=404; I=32.1; J=28.9; K=531; L=506; M=35.4; N=33.7},
    @{row=3; E=12; F=18; G=354; H=361; I=29.5; J=30.1; K=518; L=566; M=28.8; N=31.4},
    @{row=4; E=12; F=16; G=456; H=314; I=38; J=26.2; K=497; L=444; M=31.1; N=27.8},
    @{row=5; E=17; F=15; G=567; H=437; I=33.4; J=25.7; K=510; L=433; M=34; N=28.9},
    @{row=6; E=14; F=16; G=419; H=490; I=29.9; J=35; K=438; L=614; M=27.4; N=38.4},
    @{row=7; E=19; F=12; G=608; H=595; I=32; J=31.3; K=348; L=402; M=29; N=33.5},
    @{row=8; E=15; F=14; G=499; H=397; I=33.3; J=26.5; K=385; L=436; M=27.5; N=31.1},
    @{row=9; E=18; F=12; G=661; H=483; I=36.7; J=26.8; K=412; L=335; M=34.3; N=27.9},
    @{row=10; E=12; F=18; G=342; H=416; I=28.5; J=34.7; K=495; L=657; M=27.5; N=36.5},
    @{row=11; E=13; F=17; G=353; H=467; I=27.2; J=35.9; K=404; L=648; M=23.8; N=38.1},
    @{row=12; E=16; F=17; G=524; H=421; I=32.8; J=26.3; K=524; L=420; M=30.8; N=24.7},
    @{row=13; E=16; F=14; G=607; H=398; I=37.9; J=24.9; K=411; L=415; M=29.4; N=29.6},
    @{row=14; E=19; F=13; G=601; H=667; I=31.6; J=35.1; K=374; L=475; M=28.8; N=36.5},
    @{row=15; E=15; F=14; G=488; H=496; I=32.5; J=33.1; K=414; L=460; M=29.6; N=32.9},
    @{row=16; E=12; F=18; G=327; H=347; I=27.2; J=28.9; K=512; L=539; M=28.4; N=29.9},
    @{row=17; E=16; F=15; G=453; H=376; I=28.3; J=23.5; K=493; L=402; M=32.9; N=26.8},
    @{row=18; E=12; F=18; G=328; H=434; I=27.3; J=36.2; K=472; L=607; M=26.2; N=33.7},
    @{row=19; E=18; F=12; G=606; H=511; I=33.7; J=28.4; K=393; L=408; M=32.8; N=34},
    @{row=20; E=17; F=15; G=565; H=530; I=33.2; J=31.2; K=507; L=472; M=33.8; N=31.5},
    @{row=21; E=16; F=15; G=593; H=468; I=37.1; J=29.2; K=510; L=476; M=34; N=31.7},
    @{row=22; E=16; F=13; G=472; H=512; I=29.5; J=32; K=348; L=483; M=26.8; N=37.2},
    @{row=23; E=14; F=16; G=351; H=386; I=25.1; J=27.6; K=414; L=425; M=25.9; N=26.6}
)

foreach ($d in $shotsHa) {
    $ws2.Cells.Item($d.row, 4).Value = $asOfUtc
    $ws2.Cells.Item($d.row, 5).Value = $d.E
    $ws2.Cells.Item($d.row, 6).Value = $d.F
    $ws2.Cells.Item($d.row, 7).Value = $d.G
    $ws2.Cells.Item($d.row, 8).Value = $d.H
    $ws2.Cells.Item($d.row, 9).Value = $d.I
    $ws2.Cells.Item($d.row, 10).Value = $d.J
    $ws2.Cells.Item($d.row, 11).Value = $d.K
    $ws2.Cells.Item($d.row, 12).Value = $d.L
    $ws2.Cells.Item($d.row, 13).Value = $d.M
    $ws2.Cells.Item($d.row, 14).Value = $d.N
}

# ---------------------------------------------------------------------
# 3. Shots_Summary - per-team total shots-on-goal, as of 2025-11-28
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Shots_Summary")

$shotsSummary = @(
    @{row=2; E=29; F=980; G=910; H=33.8; I=31.4},
    @{row=3; E=30; F=872; G=927; H=29.1; I=30.9},
    @{row=4; E=28; F=953; G=758; H=34; I=27.1},
    @{row=5; E=32; F=1077; G=870; H=33.7; I=27.2},
    @{row=6; E=30; F=857; G=1104; H=28.6; I=36.8},
    @{row=7; E=31; F=956; G=997; H=30.8; I=32.2},
    @{row=8; E=29; F=884; G=833; H=30.5; I=28.7},
    @{row=9; E=30; F=1073; G=818; H=35.8; I=27.3},
    @{row=10; E=30; F=837; G=1073; H=27.9; I=35.8},
    @{row=11; E=30; F=757; G=1115; H=25.2; I=37.2},
    @{row=12; E=33; F=1048; G=841; H=31.8; I=25.5},
    @{row=13; E=30; F=1018; G=813; H=33.9; I=27.1},
    @{row=14; E=32; F=975; G=1142; H=30.5; I=35.7},
    @{row=15; E=29; F=902; G=956; H=31.1; I=33},
    @{row=16; E=30; F=839; G=886; H=28; I=29.5},
    @{row=17; E=31; F=946; G=778; H=30.5; I=25.1},
    @{row=18; E=30; F=800; G=1041; H=26.7; I=34.7},
    @{row=19; E=30; F=999; G=919; H=33.3; I=30.6},
    @{row=20; E=32; F=1072; G=1002; H=33.5; I=31.3},
    @{row=21; E=31; F=1103; G=944; H=35.6; I=30.5},
    @{row=22; E=29; F=820; G=995; H=28.3; I=34.3},
    @{row=23; E=30; F=765; G=811; H=25.5; I=27}
)

foreach ($d in $shotsSummary) {
    $ws3.Cells.Item($d.row, 4).Value = $asOfUtc
    $ws3.Cells.Item($d.row, 5).Value = $d.E
    $ws3.Cells.Item($d.row, 6).Value = $d.F
    $ws3.Cells.Item($d.row, 7).Value = $d.G
    $ws3.Cells.Item($d.row, 8).Value = $d.H
    $ws3.Cells.Item($d.row, 9).Value = $d.I
}

# ---------------------------------------------------------------------
# 4. Meta_ext - bump as_of_utc and build_version
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Meta_ext")
$ws4.Range("B2").Value = $asOfUtc
$ws4.Range("D2").Value = 19

Write-Output "applied 2025-11-28 KHL update"
